# Task 460, amount_applied_for and flag
# Adds two new columns (M: amount_applied_for_flag, N: amount_applied_for)
# to the "invalid" sheet (2nd sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Header row
$ws.Range("M1").Value = "amount_applied_for_flag"
$ws.Range("N1").Value = "amount_applied_for"

# Data rows: amount_applied_for_flag (M) / amount_applied_for (N)
$values = @(
    @(2, 900, 0),
    @(3, 988, $null),
    @(4, 999, $null),
    @(5, $null, 1),
    @(6, 988, 20),
    @(7, 999, 20),
    @(8, 900, 1),
    @(9, 900, 30),
    @(10, 900, $null),
    @(11, $null, 20)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Range("M$r").Value = $row[1]
    $ws.Range("N$r").Value = $row[2]
}

# New cells should carry the default "Normal" style (no explicit style
# index), matching the rest of the freshly-added M:N block.
$ws.Range("M1:N11").Style = "Normal"

# Selection now spans the newly added columns.
$ws.Range("M1:N11").Select() | Out-Null
